$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 1249.1111
$ws.Range("I40").Value = 1280.875
$ws.Range("K40").Value = 1280.875
$ws.Range("M40").Value = -1105.875
# Row 98
$ws.Range("H98").Value = 1810.9
$ws.Range("I98").Value = 800.8570999999999
$ws.Range("J98").Value = 4167.6665
$ws.Range("K98").Value = 800.8570999999999
$ws.Range("L98").Value = 4167.6665
$ws.Range("M98").Value = 697.1429000000001
$ws.Range("N98").Value = -7163.6665
# Row 106
$ws.Range("H106").Value = 9197.4
$ws.Range("I106").Value = 10747.5
$ws.Range("K106").Value = 10747.5
$ws.Range("M106").Value = -10116.5
# Row 114
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
# Row 116
$ws.Range("H116").Value = 4500
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()
# Row 122
$ws.Range("H122").Value = 1810.9
$ws.Range("I122").Value = 800.8570999999999
$ws.Range("J122").Value = 4167.6665
$ws.Range("K122").Value = 2402.5713
$ws.Range("L122").Value = 12502.9995
$ws.Range("M122").Value = 47.42870000000039
$ws.Range("N122").Value = -17402.9995
# Row 127
$ws.Range("H127").Value = 1017.5
$ws.Range("I127").Value = 1053.3334
$ws.Range("J127").Value = 910
$ws.Range("K127").Value = 3160.0002
$ws.Range("L127").Value = 2730
$ws.Range("M127").Value = 1799.9998
$ws.Range("N127").Value = -12650
# Row 132
$ws.Range("H132").Value = 3733.0454
$ws.Range("I132").Value = 3933.85
$ws.Range("J132").Value = 1725
$ws.Range("K132").Value = 11801.55
$ws.Range("L132").Value = 5175
$ws.Range("M132").Value = -9271.549999999999
$ws.Range("N132").Value = -10235
# Row 141
$ws.Range("H141").Value = 2135.6428
$ws.Range("I141").Value = 2135.6428
$ws.Range("K141").Value = 6406.928400000001
$ws.Range("M141").Value = -1226.928400000001

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 8410.581
$ws.Range("I32").Value = 8357.634
$ws.Range("K32").Value = 8357.634
$ws.Range("M32").Value = -8070.634
# Row 45
$ws.Range("H45").Value = 2735.1667
$ws.Range("I45").Value = 2735.1667
$ws.Range("K45").Value = 2735.1667
$ws.Range("M45").Value = -2358.1667
# Row 61
$ws.Range("H61").Value = 3910.1667
$ws.Range("I61").Value = 3538.3635
$ws.Range("K61").Value = 3538.3635
$ws.Range("M61").Value = -3326.3635
# Row 74
$ws.Range("H74").Value = 2051.4375
$ws.Range("I74").Value = 2021.5333
$ws.Range("J74").Value = 2500
$ws.Range("K74").Value = 2021.5333
$ws.Range("L74").Value = 2500
$ws.Range("M74").Value = -1147.5333
$ws.Range("N74").Value = -4248
# Row 77
$ws.Range("H77").Value = 2051.4375
$ws.Range("I77").Value = 2021.5333
$ws.Range("J77").Value = 2500
$ws.Range("K77").Value = 10107.6665
$ws.Range("L77").Value = 12500
$ws.Range("M77").Value = -5739.666499999999
$ws.Range("N77").Value = -21236
# Row 136
$ws.Range("H136").Value = 3910.1667
$ws.Range("I136").Value = 3538.3635
$ws.Range("K136").Value = 10615.0905
$ws.Range("M136").Value = -8065.0905

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 2242
$ws.Range("J20").Value = 4101.5
$ws.Range("L20").Value = 4101.5
$ws.Range("N20").Value = -4595.5
# Row 94
$ws.Range("H94").Value = 376.33334
$ws.Range("I94").Value = 388.0909
$ws.Range("K94").Value = 388.0909
$ws.Range("M94").Value = 62.90910000000002
# Row 105
$ws.Range("H105").Value = 4462.517
$ws.Range("I105").Value = 4404.6665
$ws.Range("J105").Value = 4557.1816
$ws.Range("K105").Value = 4404.6665
$ws.Range("L105").Value = 4557.1816
$ws.Range("M105").Value = -2657.6665
$ws.Range("N105").Value = -8051.1816
# Row 134
$ws.Range("H134").Value = 6270.4287
$ws.Range("I134").Value = 7117.3
$ws.Range("K134").Value = 21351.9
$ws.Range("M134").Value = -18816.9

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 6202.4
$ws.Range("I58").Value = 5253
$ws.Range("K58").Value = 5253
$ws.Range("M58").Value = -5050
# Row 62
$ws.Range("H62").Value = 7499.5
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 7499.5
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 7499.5
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -8747.5
# Row 65
$ws.Range("H65").Value = 7499.5
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 7499.5
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 37497.5
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -43737.5
# Row 132
$ws.Range("H132").Value = 3989.5
$ws.Range("I132").Value = 3989.5
$ws.Range("K132").Value = 11968.5
$ws.Range("M132").Value = -9438.5
# Row 136
$ws.Range("H136").Value = 6202.4
$ws.Range("I136").Value = 5253
$ws.Range("K136").Value = 15759
$ws.Range("M136").Value = -13209

$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()
# Row 34
$ws.Range("H34").Value = 56562.65
$ws.Range("I34").Value = 30000
$ws.Range("J34").Value = 57960.684
$ws.Range("K34").Value = 90000
$ws.Range("L34").Value = 173882.052
$ws.Range("M34").Value = -89916
$ws.Range("N34").Value = -174050.052
# Row 131
$ws.Range("H131").Value = 1773
$ws.Range("J131").Value = 2033
$ws.Range("L131").Value = 6099
$ws.Range("N131").Value = -16179
# Row 140
$ws.Range("H140").Value = 1002439.4
$ws.Range("I140").Value = 1002439.4
$ws.Range("K140").Value = 3007318.2
$ws.Range("M140").Value = -3002138.2

$ws = $wb.Worksheets.Item("GSM")
# Row 93
$ws.Range("H93").Value = 68999.5
$ws.Range("J93").Value = 68999.5
$ws.Range("L93").Value = 68999.5
$ws.Range("N93").Value = -72743.5

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 1763.6923
$ws.Range("I46").Value = 1478.75
$ws.Range("J46").Value = 2219.6
$ws.Range("K46").Value = 1478.75
$ws.Range("L46").Value = 2219.6
$ws.Range("M46").Value = -1290.75
$ws.Range("N46").Value = -2595.6
# Row 132
$ws.Range("H132").Value = 6962.25
$ws.Range("I132").Value = 6962.25
$ws.Range("K132").Value = 20886.75
$ws.Range("M132").Value = -18356.75
# Row 136
$ws.Range("H136").Value = 2999.1428
$ws.Range("I136").Value = 2665.6667
$ws.Range("K136").Value = 7997.000100000001
$ws.Range("M136").Value = -5447.000100000001

$ws = $wb.Worksheets.Item("WVR")
# Row 45
$ws.Range("H45").Value = 95312.5
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 95312.5
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 95312.5
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -96294.5
# Row 122
$ws.Range("H122").Value = 3998.7144
$ws.Range("I122").Value = 3998.7144
$ws.Range("K122").Value = 11996.1432
$ws.Range("M122").Value = -9546.143199999999
# Row 132
$ws.Range("H132").Value = 1754.8422
$ws.Range("I132").Value = 1794.3243
$ws.Range("K132").Value = 5382.9729
$ws.Range("M132").Value = -2852.9729
# Row 136
$ws.Range("H136").Value = 5197.3477
$ws.Range("I136").Value = 4681.9287
$ws.Range("J136").Value = 5999.1113
$ws.Range("K136").Value = 14045.7861
$ws.Range("L136").Value = 17997.3339
$ws.Range("M136").Value = -11495.7861
$ws.Range("N136").Value = -23097.3339
